$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values (text strings) for the daily KPI pagination.
# Old dates 2023-09-18 .. 2023-09-23 (rows 2-49, in blocks of 8 rows)
# become 2023-09-25 .. 2023-09-30 (shifted by one week).
$dateMap = @{
    "2023-09-18" = "2023-09-25"
    "2023-09-19" = "2023-09-26"
    "2023-09-20" = "2023-09-27"
    "2023-09-21" = "2023-09-28"
    "2023-09-22" = "2023-09-29"
    "2023-09-23" = "2023-09-30"
}

for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    if ($dateMap.ContainsKey($old)) {
        $cell.Value = $dateMap[$old]
    }
}

# Update the sheet view scroll/selection state.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B35").Select()
